$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`nBinance`n✅ 1000 Bs = 9.28 = 39034.8 pesos`n✅ 39034.8 pesos = 9.26 = 969.71 Bs`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate figures in N10/O10 and N12/O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 107.75
$ws2.Range("O10").Value = 4206
$ws2.Range("N12").Value = 4215
$ws2.Range("O12").Value = 104.71
